$wb = $excel.ActiveWorkbook

# Trade #4 - new row 5 on both the "All Trades" and "base_strategy" sheets.
$rowValues = @(4, "2026-02-16", "22:57:37", "base_strategy", "DOWN", 0.5, "", "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", "", 0)
$emptyTextCols = @(7, 16)   # G (Exit Price) and P (Exit Reason) stay as empty text, like existing OPEN rows

foreach ($sheetName in @("All Trades", "base_strategy")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item(5, $col)
        if ($col -eq 2) {
            # Column B ("Date") holds a pure ISO date string like earlier rows;
            # force text so Excel doesn't auto-convert it to a date serial.
            $cell.NumberFormat = "@"
            $cell.Value = $rowValues[$i]
            $cell.Style = "Normal"
        } elseif ($emptyTextCols -contains $col) {
            # Keep these as an explicit empty text value (matching the other
            # OPEN-status rows' empty Exit Price / Exit Reason cells) instead
            # of leaving the cell completely blank.
            $cell.Value = "'"
            $cell.Style = "Normal"
        } else {
            $cell.Value = $rowValues[$i]
        }
    }
}
